$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new values in row 3
$ws.Range("A3").Value = "2,3,4"
$ws.Range("B3").Value = 0.5
$ws.Range("C3").Value = "Checking outlliers`nBivariate and multivariate analysis, KNN, SelectKBest, GridSearchCV"

# Update the selection to C5
$ws.Range("C5").Select()
